# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" worksheets, which contain duplicated data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1793
    $ws.Range("F5").Value = 1110
    $ws.Range("F6").Value = 995
    $ws.Range("F8").Value = 5905
}
